$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (generated from diff)
$ws.Range("D2").Value = "62.985.36"
$ws.Range("E2").Value = "  -2.15%  "
$ws.Range("D3").Value = "2.675.58"
$ws.Range("E3").Value = "  -3.67%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.12"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -5.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.88"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.44%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("E8").Value = "  -2.76%  "
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("E11").Value = "  -5.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.08"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -14.07%  "
$ws.Range("D13").Value = "3.144.77"
$ws.Range("E13").Value = "  -3.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.04"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.14%  "
$ws.Range("D15").Value = "62.874.56"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("E16").Value = "  -5.80%  "
$ws.Range("D17").Value = "2.677.69"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.90"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("E19").Value = "  -6.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "340.89"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -5.79%  "
$ws.Range("E21").Value = "  -6.59%  "
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("E23").Value = "  -5.74%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.41"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.08"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.70%  "
$ws.Range("D28").Value = "0.0₃0846"
$ws.Range("E28").Value = "  -8.08%  "
$ws.Range("E29").Value = "  -3.55%  "
$ws.Range("E30").Value = "  -3.04%  "
$ws.Range("E31").Value = "  -5.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.71"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.67%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.998"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -5.01%  "
$ws.Range("E35").Value = "  -3.95%  "
$ws.Range("E36").Value = "  -6.41%  "
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "337.72"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.97%  "
$ws.Range("E39").Value = "  -7.61%  "
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.05"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.93%  "
$ws.Range("E42").Value = "  -7.29%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.30%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.14"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -7.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.616"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -2.99%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("E47").Value = "  -6.53%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "11.04"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.10%  "
$ws.Range("E49").Value = "  -4.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.37"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -6.83%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0238"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -6.65%  "
